# Edit workbook: update price rows 70-75 and append new rows 76-79
# for "Hortaliza, Agrícola del Norte S.A. de Arica - Sandia" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 70: update date, quality stays Primera, update volume/price columns ---
$ws.Range("D70").Value = 45212
$ws.Range("J70").Value = 300
$ws.Range("K70").Value = 500
$ws.Range("L70").Value = 550
$ws.Range("M70").Value = 533
$ws.Range("P70").Value = 533

# --- Row 71: update date, quality Tercera -> Segunda, update price columns ---
$ws.Range("D71").Value = 45212
$ws.Range("I71").Value = "Segunda"
$ws.Range("J71").Value = 340
$ws.Range("K71").Value = 500
$ws.Range("L71").Value = 550
$ws.Range("M71").Value = 528
$ws.Range("P71").Value = 528

# --- Row 72: update date, quality Primera -> Tercera, update price columns ---
$ws.Range("D72").Value = 45212
$ws.Range("I72").Value = "Tercera"
$ws.Range("J72").Value = 220
$ws.Range("K72").Value = 500
$ws.Range("L72").Value = 550
$ws.Range("M72").Value = 523
$ws.Range("P72").Value = 523

# --- Row 73: update date, keep quality Primera, update price columns ---
$ws.Range("D73").Value = 44214
$ws.Range("K73").Value = 400
$ws.Range("L73").Value = 450
$ws.Range("M73").Value = 425
$ws.Range("P73").Value = 425

# --- Row 74: update date, keep quality Tercera, update volume/price columns, origin ---
$ws.Range("D74").Value = 44900
$ws.Range("J74").Value = 600
$ws.Range("K74").Value = 480
$ws.Range("L74").Value = 500
$ws.Range("M74").Value = 490
$ws.Range("O74").Value = "Perú"
$ws.Range("P74").Value = 490

# --- Row 75: update date, keep quality Primera, update price columns ---
$ws.Range("D75").Value = 44251
$ws.Range("K75").Value = 250
$ws.Range("L75").Value = 280
$ws.Range("M75").Value = 265
$ws.Range("P75").Value = 265

# --- Row 76: new row appended ---
$ws.Range("A76").Value = 1
$ws.Range("B76").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C76").Value = "Arica y Parinacota"
$ws.Range("D76").Value = 44609
$ws.Range("E76").Value = 15
$ws.Range("F76").Value = 100112028
$ws.Range("G76").Value = "Sandia"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 1200
$ws.Range("K76").Value = 280
$ws.Range("L76").Value = 300
$ws.Range("M76").Value = 290
$ws.Range("N76").Value = "$/kilo (volumen en unidades)"
$ws.Range("O76").Value = "Perú"
$ws.Range("P76").Value = 290
$ws.Range("Q76").Value = 1
$ws.Range("R76").Value = "Hortaliza"

# --- Row 77: new row appended ---
$ws.Range("A77").Value = 1
$ws.Range("B77").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C77").Value = "Arica y Parinacota"
$ws.Range("D77").Value = 44586
$ws.Range("E77").Value = 15
$ws.Range("F77").Value = 100112028
$ws.Range("G77").Value = "Sandia"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Tercera"
$ws.Range("J77").Value = 500
$ws.Range("K77").Value = 330
$ws.Range("L77").Value = 350
$ws.Range("M77").Value = 340
$ws.Range("N77").Value = "$/kilo (volumen en unidades)"
$ws.Range("O77").Value = "Región de Arica y Parinacota"
$ws.Range("P77").Value = 340
$ws.Range("Q77").Value = 1
$ws.Range("R77").Value = "Hortaliza"

# --- Row 78: new row appended ---
$ws.Range("A78").Value = 1
$ws.Range("B78").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C78").Value = "Arica y Parinacota"
$ws.Range("D78").Value = 44243
$ws.Range("E78").Value = 15
$ws.Range("F78").Value = 100112028
$ws.Range("G78").Value = "Sandia"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 1200
$ws.Range("K78").Value = 300
$ws.Range("L78").Value = 320
$ws.Range("M78").Value = 310
$ws.Range("N78").Value = "$/kilo (volumen en unidades)"
$ws.Range("O78").Value = "Perú"
$ws.Range("P78").Value = 310
$ws.Range("Q78").Value = 1
$ws.Range("R78").Value = "Hortaliza"

# --- Row 79: new row appended ---
$ws.Range("A79").Value = 1
$ws.Range("B79").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C79").Value = "Arica y Parinacota"
$ws.Range("D79").Value = 44243
$ws.Range("E79").Value = 15
$ws.Range("F79").Value = 100112028
$ws.Range("G79").Value = "Sandia"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Segunda"
$ws.Range("J79").Value = 800
$ws.Range("K79").Value = 300
$ws.Range("L79").Value = 320
$ws.Range("M79").Value = 310
$ws.Range("N79").Value = "$/kilo (volumen en unidades)"
$ws.Range("O79").Value = "Perú"
$ws.Range("P79").Value = 310
$ws.Range("Q79").Value = 1
$ws.Range("R79").Value = "Hortaliza"

# Ensure date cells use the same date number format as the rest of column D
$ws.Range("D76").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D77").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D78").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D79").NumberFormat = "YYYY-MM-DD HH:MM:SS"
